$wb = $excel.ActiveWorkbook

# ==========================================================================
# Refresh the "想去人数" (want-to-go count) / "最低票价" (min price) figures
# across the four sheets, mirroring the upstream bilibili-show scrape that
# produced this gh-pages snapshot.
# ==========================================================================

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 154
$ws.Range("F7").Value = 5522
$ws.Range("F10").Value = 3789
$ws.Range("F11").Value = 67
$ws.Range("F18").Value = 99
$ws.Range("F19").Value = 133
$ws.Range("F20").Value = 173
$ws.Range("G21").Value = 69
$ws.Range("F23").Value = 5193
$ws.Range("F25").Value = 2062
$ws.Range("F26").Value = 126
$ws.Range("F27").Value = 336
$ws.Range("F28").Value = 7698
$ws.Range("F29").Value = 31
$ws.Range("F31").Value = 2182
$ws.Range("F32").Value = 2151
$ws.Range("F33").Value = 1325
$ws.Range("F34").Value = 154
$ws.Range("F35").Value = 1179
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 257
$ws.Range("F40").Value = 241
$ws.Range("F42").Value = 1174
$ws.Range("F43").Value = 1171
$ws.Range("F45").Value = 1318
$ws.Range("F46").Value = 2024
$ws.Range("F47").Value = 116
$ws.Range("F48").Value = 207
$ws.Range("F49").Value = 1210

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 117

# New event row appended at the bottom (row 20).
# Copy the format of each cell in the row above first (single-cell ranges,
# NOT whole-row copies - copying an entire Rows.Item would paste formatting
# across all 16384 columns and bloat the sheet) so the new row's styling
# (bold/centered/bordered A-column, plain text B..I columns) matches the
# rest of the table, then fill in the values.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(20, 1).Value = 19

# Column B holds a plain "yyyy-mm-dd" text label, not a real date - force
# text so Excel doesn't reinterpret it as a date serial number.
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = "2024-10-26"
# Restore the default (General) number format on B20 now that the text
# value is safely stored, so its styling matches the other date cells.
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item(20, 3).Value = "北京·变形金刚音乐会40周年特别版"
$ws.Cells.Item(20, 4).Value = "中关村南大街33号国家图书馆北门 国图艺术中心音乐厅"
$ws.Cells.Item(20, 5).Value = "2024.10.26 19:30-10.26 21:30"
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 171
$ws.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89970"
$ws.Cells.Item(20, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/TwvRQI041722150343639.jpeg"

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 548
$ws.Range("F3").Value = 722

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 154
$ws.Range("F6").Value = 548
$ws.Range("F7").Value = 722
$ws.Range("F9").Value = 5522
$ws.Range("F10").Value = 3789
$ws.Range("F11").Value = 67
$ws.Range("F17").Value = 99
$ws.Range("F20").Value = 173
$ws.Range("G22").Value = 69
$ws.Range("F24").Value = 5193
$ws.Range("F26").Value = 2062
$ws.Range("F27").Value = 126
$ws.Range("F28").Value = 336
$ws.Range("F29").Value = 7698
$ws.Range("F30").Value = 31
$ws.Range("F32").Value = 2182
$ws.Range("F33").Value = 2151
$ws.Range("F34").Value = 1325
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 1179
$ws.Range("F37").Value = 257
$ws.Range("F38").Value = 241
$ws.Range("F40").Value = 1174
$ws.Range("F41").Value = 1171
$ws.Range("F43").Value = 1318
$ws.Range("F45").Value = 2024
$ws.Range("F46").Value = 116
$ws.Range("F48").Value = 207
$ws.Range("F49").Value = 1210
